# Lab 2 Report - fill in the "Without Correction" Delta X / Delta Y data
# (columns C/D) and the first data row of "With Correction" (columns G/H)
# that was missing, so the Standard Deviation formulas (J5:M5) compute
# real numbers instead of #DIV/0! errors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value  = 0
$ws.Range("D4").Value  = 0.1
$ws.Range("G4").Value  = 0
$ws.Range("H4").Value  = 0

$ws.Range("C5").Value  = -0.3
$ws.Range("D5").Value  = 0.5

$ws.Range("C6").Value  = -0.5
$ws.Range("D6").Value  = 0.5

$ws.Range("C7").Value  = -0.1
$ws.Range("D7").Value  = 0

$ws.Range("C8").Value  = 0.1
$ws.Range("D8").Value  = 0

$ws.Range("C9").Value  = 0
$ws.Range("D9").Value  = -0.1

$ws.Range("C10").Value = -0.2
$ws.Range("D10").Value = 0

$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.1

$ws.Range("C12").Value = 0.1
$ws.Range("D12").Value = -0.3

$ws.Range("C13").Value = -0.1
$ws.Range("D13").Value = -0.3

# Recalculate workbook so the STDEV.P formulas in J5:M5 re-evaluate with
# the newly entered data instead of remaining cached #DIV/0! errors.
$excel.CalculateFull()

# Leave the selection where the author ended up after entering the data.
$ws.Range("E13").Select()

$wb.Save()
